$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Update price values in column D (ARA.PLANA section, rows 33-47)
$ws.Range("D33").Value = 5030.479
$ws.Range("D34").Value = 3779.603
$ws.Range("D35").Value = 3261.594
$ws.Range("D36").Value = 2929.043
$ws.Range("D37").Value = 2929.043
$ws.Range("D38").Value = 2470.491
$ws.Range("D39").Value = 2470.491
$ws.Range("D40").Value = 2470.491
$ws.Range("D41").Value = 2470.491
$ws.Range("D42").Value = 2470.491
$ws.Range("D43").Value = 2470.491
$ws.Range("D44").Value = 2470.491
$ws.Range("D45").Value = 2948.224
$ws.Range("D46").Value = 2948.224
$ws.Range("D47").Value = 2948.224

# Update price values in column D (ARANDELA CHAPISTA section, rows 51-54)
$ws.Range("D51").Value = 4297.629
$ws.Range("D52").Value = 3907.523
$ws.Range("D53").Value = 3261.594
$ws.Range("D54").Value = 3261.594
